$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New candidate "徐媛媛" is appended as row 13, following the exact same
# layout/formatting as row 12 (the previous last row). Add the mailto
# hyperlink first, then copy row 12's formatting+values down onto row 13
# (this keeps B13's style identical to the other e-mail cells instead of
# picking up the extra "applyFont" that Hyperlinks.Add stamps in place).
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:329304319@qq.com")
$ws.Range("A12:W12").Copy($ws.Range("A13:W13"))
$ws.Rows.Item(13).RowHeight = 40.5

# Candidate info
$ws.Range("A13").Value = "徐媛媛"
$ws.Range("B13").Value = "329304319@qq.com"
$ws.Range("D13").Value = 43590
$ws.Range("E13").Value = 17765101962
$ws.Range("F13").Value = "本科"
$ws.Range("H13").Value = 10
$ws.Range("J13").Value = '已收到，酷家乐效果图网络连接：https://yun.kujiale.com/design/3FO41MHL5TVQ/show?fromqrcode=true&from=panoMp'
$ws.Range("L13").Value = "N"
$ws.Range("M13").ClearContents()
